# Swap the presentation's theme colour scheme from the "Integral" palette
# to the stock "Office Theme" palette (ppt/theme/theme1.xml), mirroring the
# canonical OOXML edit where theme1.xml/theme2.xml swap contents.
#
# PowerPoint's COM colour scheme is addressed as:
#   Master.Theme.ThemeColorScheme.Item(index).RGB
# with index order: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3
#                    8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
# RGB is a COM "long" packed as 0x00BBGGRR (blue high byte, red low byte).

$p = $ppt.ActivePresentation

function ToComRgb($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

# Target "Office Theme" palette (hex RRGGBB), in clrScheme slot order.
$officeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 0; $i -lt $officeHex.Count; $i++) {
    $hex = $officeHex[$i]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Item($i + 1).RGB = ToComRgb $r $g $b
}

Write-Output "Theme colour scheme updated to Office Theme palette."
